$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(18, 8).Value = 166833570
$ws.Cells.Item(18, 9).Value = 166833570
$ws.Cells.Item(18, 11).Value = 166833570
$ws.Cells.Item(18, 13).Value = -166833286
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(5, 8).Value = 79.09090999999999
$ws.Cells.Item(5, 9).Value = 56.25
$ws.Cells.Item(5, 10).Value = 140
$ws.Cells.Item(5, 11).Value = 56.25
$ws.Cells.Item(5, 12).Value = 140
$ws.Cells.Item(5, 13).Value = 55.75
$ws.Cells.Item(5, 14).Value = -364
$ws.Cells.Item(9, 8).Value = 19666.666
$ws.Cells.Item(9, 10).Value = 19666.666
$ws.Cells.Item(9, 12).Value = 19666.666
$ws.Cells.Item(9, 14).Value = -20006.666
$ws.Cells.Item(20, 8).Value = 19666.666
$ws.Cells.Item(20, 10).Value = 19666.666
$ws.Cells.Item(20, 12).Value = 19666.666
$ws.Cells.Item(20, 14).Value = -20206.666
$ws.Cells.Item(74, 8).Value = 1164.0227
$ws.Cells.Item(74, 9).Value = 1184.081
$ws.Cells.Item(74, 11).Value = 1184.081
$ws.Cells.Item(74, 13).Value = -310.0809999999999
$ws.Cells.Item(77, 8).Value = 1164.0227
$ws.Cells.Item(77, 9).Value = 1184.081
$ws.Cells.Item(77, 11).Value = 5920.405
$ws.Cells.Item(77, 13).Value = -1552.405
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(4, 8).Value = 79.09090999999999
$ws.Cells.Item(4, 9).Value = 56.25
$ws.Cells.Item(4, 10).Value = 140
$ws.Cells.Item(4, 11).Value = 56.25
$ws.Cells.Item(4, 12).Value = 140
$ws.Cells.Item(4, 13).Value = 58.75
$ws.Cells.Item(4, 14).Value = -370
$ws.Cells.Item(10, 8).Value = 176.25
$ws.Cells.Item(10, 9).Value = 176.25
$ws.Cells.Item(10, 10).Value = 0
$ws.Cells.Item(10, 11).Value = 176.25
$ws.Cells.Item(10, 12).Value = 0
$ws.Cells.Item(10, 13).Value = -36.25
$ws.Cells.Item(10, 14).Value = ""
$ws.Cells.Item(18, 8).Value = 50000
$ws.Cells.Item(18, 10).Value = 50000
$ws.Cells.Item(18, 12).Value = 50000
$ws.Cells.Item(18, 14).Value = -51058
$ws.Cells.Item(20, 8).Value = 1940.697
$ws.Cells.Item(20, 9).Value = 1350.8
$ws.Cells.Item(20, 11).Value = 1350.8
$ws.Cells.Item(20, 13).Value = -1103.8
$ws.Cells.Item(22, 8).Value = 0
$ws.Cells.Item(22, 9).Value = 0
$ws.Cells.Item(22, 11).Value = 0
$ws.Cells.Item(22, 13).Value = ""
$ws.Cells.Item(134, 8).Value = 1728.9375
$ws.Cells.Item(134, 9).Value = 1766.4231
$ws.Cells.Item(134, 10).Value = 1566.5
$ws.Cells.Item(134, 11).Value = 5299.2693
$ws.Cells.Item(134, 12).Value = 4699.5
$ws.Cells.Item(134, 13).Value = -2764.2693
$ws.Cells.Item(134, 14).Value = -9769.5
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 1953.0408
$ws.Cells.Item(31, 9).Value = 1398.8948
$ws.Cells.Item(31, 11).Value = 1398.8948
$ws.Cells.Item(31, 13).Value = -1103.8948
$ws.Cells.Item(34, 8).Value = 1953.0408
$ws.Cells.Item(34, 9).Value = 1398.8948
$ws.Cells.Item(34, 11).Value = 1398.8948
$ws.Cells.Item(34, 13).Value = -1196.8948
$ws.Cells.Item(58, 8).Value = 823.0909
$ws.Cells.Item(58, 9).Value = 695.4902
$ws.Cells.Item(58, 10).Value = 2450
$ws.Cells.Item(58, 11).Value = 695.4902
$ws.Cells.Item(58, 12).Value = 2450
$ws.Cells.Item(58, 13).Value = -492.4902
$ws.Cells.Item(58, 14).Value = -2856
$ws.Cells.Item(136, 8).Value = 823.0909
$ws.Cells.Item(136, 9).Value = 695.4902
$ws.Cells.Item(136, 10).Value = 2450
$ws.Cells.Item(136, 11).Value = 2086.4706
$ws.Cells.Item(136, 12).Value = 7350
$ws.Cells.Item(136, 13).Value = 463.5294000000004
$ws.Cells.Item(136, 14).Value = -12450
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(40, 8).Value = 396.15384
$ws.Cells.Item(40, 9).Value = 233.33333
$ws.Cells.Item(40, 10).Value = 445
$ws.Cells.Item(40, 11).Value = 933.33332
$ws.Cells.Item(40, 12).Value = 1780
$ws.Cells.Item(40, 13).Value = -864.33332
$ws.Cells.Item(40, 14).Value = -1918
$ws.Cells.Item(62, 8).Value = 4132.2856
$ws.Cells.Item(62, 9).Value = 1112
$ws.Cells.Item(62, 10).Value = 4635.6665
$ws.Cells.Item(62, 11).Value = 3336
$ws.Cells.Item(62, 12).Value = 13906.9995
$ws.Cells.Item(62, 13).Value = -2650
$ws.Cells.Item(62, 14).Value = -15278.9995
$ws.Cells.Item(65, 8).Value = 4132.2856
$ws.Cells.Item(65, 9).Value = 1112
$ws.Cells.Item(65, 10).Value = 4635.6665
$ws.Cells.Item(65, 11).Value = 10008
$ws.Cells.Item(65, 12).Value = 41720.9985
$ws.Cells.Item(65, 13).Value = -6576
$ws.Cells.Item(65, 14).Value = -48584.9985
$ws.Cells.Item(92, 8).Value = 800
$ws.Cells.Item(92, 9).Value = 0
$ws.Cells.Item(92, 10).Value = 800
$ws.Cells.Item(92, 11).Value = 0
$ws.Cells.Item(92, 12).Value = 2400
$ws.Cells.Item(92, 13).Value = ""
$ws.Cells.Item(92, 14).Value = -4896
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(2, 8).Value = 84.09090999999999
$ws.Cells.Item(2, 9).Value = 47.22222
$ws.Cells.Item(2, 10).Value = 250
$ws.Cells.Item(2, 11).Value = 47.22222
$ws.Cells.Item(2, 12).Value = 250
$ws.Cells.Item(2, 13).Value = 65.77778000000001
$ws.Cells.Item(2, 14).Value = -476
$ws.Cells.Item(28, 8).Value = 0
$ws.Cells.Item(28, 10).Value = 0
$ws.Cells.Item(28, 12).Value = 0
$ws.Cells.Item(28, 14).Value = ""
$ws.Cells.Item(70, 8).Value = 5947.788
$ws.Cells.Item(70, 9).Value = 6195.1113
$ws.Cells.Item(70, 10).Value = 4834.8335
$ws.Cells.Item(70, 11).Value = 6195.1113
$ws.Cells.Item(70, 12).Value = 4834.8335
$ws.Cells.Item(70, 13).Value = -5925.1113
$ws.Cells.Item(70, 14).Value = -5374.8335
$ws.Cells.Item(73, 8).Value = 5947.788
$ws.Cells.Item(73, 9).Value = 6195.1113
$ws.Cells.Item(73, 10).Value = 4834.8335
$ws.Cells.Item(73, 11).Value = 6195.1113
$ws.Cells.Item(73, 12).Value = 4834.8335
$ws.Cells.Item(73, 13).Value = -5259.1113
$ws.Cells.Item(73, 14).Value = -6706.8335
$ws.Cells.Item(132, 8).Value = 2781
$ws.Cells.Item(132, 9).Value = 2545.3447
$ws.Cells.Item(132, 10).Value = 3350.5
$ws.Cells.Item(132, 11).Value = 7636.034100000001
$ws.Cells.Item(132, 12).Value = 10051.5
$ws.Cells.Item(132, 13).Value = -5106.034100000001
$ws.Cells.Item(132, 14).Value = -15111.5
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(22, 8).Value = 784
$ws.Cells.Item(22, 9).Value = 660
$ws.Cells.Item(22, 10).Value = 1280
$ws.Cells.Item(22, 11).Value = 660
$ws.Cells.Item(22, 12).Value = 1280
$ws.Cells.Item(22, 13).Value = -365
$ws.Cells.Item(22, 14).Value = -1870
$ws.Cells.Item(27, 8).Value = 784
$ws.Cells.Item(27, 9).Value = 660
$ws.Cells.Item(27, 10).Value = 1280
$ws.Cells.Item(27, 11).Value = 660
$ws.Cells.Item(27, 12).Value = 1280
$ws.Cells.Item(27, 13).Value = -553
$ws.Cells.Item(27, 14).Value = -1494
$ws.Cells.Item(61, 8).Value = 1705.8889
$ws.Cells.Item(61, 9).Value = 1187
$ws.Cells.Item(61, 10).Value = 2121
$ws.Cells.Item(61, 11).Value = 1187
$ws.Cells.Item(61, 12).Value = 2121
$ws.Cells.Item(61, 13).Value = -985
$ws.Cells.Item(61, 14).Value = -2525
$ws.Cells.Item(93, 8).Value = 13774
$ws.Cells.Item(93, 9).Value = 20444.445
$ws.Cells.Item(93, 10).Value = 3768.3333
$ws.Cells.Item(93, 11).Value = 20444.445
$ws.Cells.Item(93, 12).Value = 3768.3333
$ws.Cells.Item(93, 13).Value = -19196.445
$ws.Cells.Item(93, 14).Value = -6264.3333
$ws.Cells.Item(113, 8).Value = 1705.8889
$ws.Cells.Item(113, 9).Value = 1187
$ws.Cells.Item(113, 10).Value = 2121
$ws.Cells.Item(113, 11).Value = 1187
$ws.Cells.Item(113, 12).Value = 2121
$ws.Cells.Item(113, 13).Value = 983
$ws.Cells.Item(113, 14).Value = -6461
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(22, 8).Value = 1000
$ws.Cells.Item(22, 9).Value = 1000
$ws.Cells.Item(22, 11).Value = 1000
$ws.Cells.Item(22, 13).Value = -707
$ws.Cells.Item(110, 8).Value = 50000
$ws.Cells.Item(110, 10).Value = 50000
$ws.Cells.Item(110, 12).Value = 50000
$ws.Cells.Item(110, 14).Value = -58180
$ws.Cells.Item(116, 8).Value = 0
$ws.Cells.Item(116, 10).Value = 0
$ws.Cells.Item(116, 12).Value = 0
$ws.Cells.Item(116, 14).Value = ""
$ws.Cells.Item(122, 8).Value = 2824.8096
$ws.Cells.Item(122, 9).Value = 3150.125
$ws.Cells.Item(122, 10).Value = 2624.6155
$ws.Cells.Item(122, 11).Value = 9450.375
$ws.Cells.Item(122, 12).Value = 7873.8465
$ws.Cells.Item(122, 13).Value = -7000.375
$ws.Cells.Item(122, 14).Value = -12773.8465
$ws.Cells.Item(132, 8).Value = 1052.8246
$ws.Cells.Item(132, 9).Value = 671.913
$ws.Cells.Item(132, 10).Value = 2015.739
$ws.Cells.Item(132, 13).Value = 514.261
